$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove extra rows 7-9 (table shrinks from 9 data+header rows to 6)
$ws.Rows("7:9").Delete()

# Update rows 2-6 with the refreshed TPM-derived values (sending cluster fixed to ECs)
# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf9"
$ws.Cells.Item(2, 3).Value = "Fgfr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.083607
$ws.Cells.Item(2, 8).Value = 3.250821
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.903185
$ws.Cells.Item(2, 14).Value = 2.709555
$ws.Cells.Item(2, 15).Value = 0.03154869388788047
$ws.Cells.Item(2, 16).Value = 0.03154869388788046
$ws.Cells.Item(2, 17).Value = 0.978697588295
$ws.Cells.Item(2, 18).Value = 8.808278294655
$ws.Cells.Item(2, 19).Value = 0.03154869388788047
$ws.Cells.Item(2, 20).Value = 0.03154869388788046

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf9"
$ws.Cells.Item(3, 3).Value = "Fgfr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.083607
$ws.Cells.Item(3, 8).Value = 3.250821
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.110372
$ws.Cells.Item(3, 14).Value = 0.331116
$ws.Cells.Item(3, 15).Value = 0.003855347953955327
$ws.Cells.Item(3, 16).Value = 0.003855347953955326
$ws.Cells.Item(3, 17).Value = 0.119599871804
$ws.Cells.Item(3, 18).Value = 1.076398846236
$ws.Cells.Item(3, 19).Value = 0.003855347953955327
$ws.Cells.Item(3, 20).Value = 0.003855347953955326

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf9"
$ws.Cells.Item(4, 3).Value = "Fgfr4"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.083607
$ws.Cells.Item(4, 8).Value = 3.250821
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.03076233333333333
$ws.Cells.Item(4, 14).Value = 0.09228699999999999
$ws.Cells.Item(4, 15).Value = 0.001074543352259254
$ws.Cells.Item(4, 16).Value = 0.001074543352259254
$ws.Cells.Item(4, 17).Value = 0.03333427973633333
$ws.Cells.Item(4, 18).Value = 0.300008517627
$ws.Cells.Item(4, 19).Value = 0.001074543352259254
$ws.Cells.Item(4, 20).Value = 0.001074543352259254

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Fgf9"
$ws.Cells.Item(5, 3).Value = "Fgfr4"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.083607
$ws.Cells.Item(5, 8).Value = 3.250821
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 27.52907633333334
$ws.Cells.Item(5, 14).Value = 82.58722900000001
$ws.Cells.Item(5, 15).Value = 0.9616041035407232
$ws.Cells.Item(5, 16).Value = 0.9616041035407231
$ws.Cells.Item(5, 17).Value = 29.83069981833433
$ws.Cells.Item(5, 18).Value = 268.476298365009
$ws.Cells.Item(5, 19).Value = 0.9616041035407232
$ws.Cells.Item(5, 20).Value = 0.9616041035407231

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Fgf9"
$ws.Cells.Item(6, 3).Value = "Fgfr4"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.083607
$ws.Cells.Item(6, 8).Value = 3.250821
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.05488933333333334
$ws.Cells.Item(6, 14).Value = 0.164668
$ws.Cells.Item(6, 15).Value = 0.001917311265181737
$ws.Cells.Item(6, 16).Value = 0.001917311265181736
$ws.Cells.Item(6, 17).Value = 0.05947846582533334
$ws.Cells.Item(6, 18).Value = 0.535306192428
$ws.Cells.Item(6, 19).Value = 0.001917311265181737
$ws.Cells.Item(6, 20).Value = 0.001917311265181736
